# Add a new "game_date" column (AU) to the "Spreads" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spreads")

# --- Header cell AU1: "game_date" with the same bold/border/centered
#     formatting used by the other header cells (e.g. AT1). ---
$headerSrc = $ws.Range("AT1")
$header = $ws.Range("AU1")
$header.Value = "game_date"
$header.Font.Bold = $headerSrc.Font.Bold
$header.HorizontalAlignment = $headerSrc.HorizontalAlignment
$header.VerticalAlignment = $headerSrc.VerticalAlignment
$header.Borders.LineStyle = $headerSrc.Borders.LineStyle

# --- Data cells AU2:AU5: plain text date strings "2025-11-12". ---
# Force text formatting first so Excel does not auto-convert the
# yyyy-mm-dd-looking string into a date serial number, then reset the
# cell style back to the default (no special formatting), matching the
# plain, unstyled data cells used elsewhere in the sheet.
$dataRange = $ws.Range("AU2:AU5")
$dataRange.NumberFormat = "@"
$dataRange.Value = "2025-11-12"
$dataRange.Style = "Normal"
